$wb = $excel.ActiveWorkbook

$wsUnitSummary = $wb.Worksheets.Item("UnitSummary")
$wsUnitMix = $wb.Worksheets.Item("UnitMix")
$wsScheme = $wb.Worksheets.Item("Scheme Summary")

# Clear out all data on the UnitSummary sheet (sheet2)
$wsUnitSummary.Cells.ClearContents()

# Clear out all data on the UnitMix sheet (sheet3)
$wsUnitMix.Cells.ClearContents()

# Reset the selection on each sheet to mimic the saved cursor positions
$wsUnitSummary.Activate()
$wsUnitSummary.Range("A1:G11").Select()

$wsUnitMix.Activate()
$wsUnitMix.Range("D19").Select()

$wsScheme.Activate()
$wsScheme.Range("A2").Select()

$wb.Save()
